$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.330446976246833
$ws.Range("C2").Value = 0.194157323651325
$ws.Range("E2").Value = 0.0936705614147576
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 1.087524282702759
$ws.Range("H2").Value = 1.056188439281073
$ws.Range("L2").Value = 0.1902377651464704
$ws.Range("M2").Value = 0.2772243007943302
$ws.Range("N2").Value = 1.593763168748349

$ws.Range("B3").Value = 1.23957629976951
$ws.Range("C3").Value = 0.1787859840349029
$ws.Range("E3").Value = 0.09411637750768043
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 1.078970107672077
$ws.Range("H3").Value = 1.058313323287422
$ws.Range("L3").Value = 0.1876884228158389
$ws.Range("M3").Value = 0.263402940416043
$ws.Range("N3").Value = 1.614770197836771

$ws.Range("B4").Value = 1.1843680810656
$ws.Range("C4").Value = 0.1692675202043574
$ws.Range("E4").Value = 0.09440676254712221
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 1.074568236877667
$ws.Range("H4").Value = 1.060231973323212
$ws.Range("L4").Value = 0.1862177759352974
$ws.Range("M4").Value = 0.2550543192679839
$ws.Range("N4").Value = 1.628317568016519

$ws.Range("B5").Value = 1.162018004448157
$ws.Range("C5").Value = 0.1653683984207532
$ws.Range("E5").Value = 0.09452929240867658
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 1.072987454399041
$ws.Range("H5").Value = 1.061167970562039
$ws.Range("L5").Value = 0.1856422979994505
$ws.Range("M5").Value = 0.2516868612507395
$ws.Range("N5").Value = 1.634001296124298

$ws.Range("B6").Value = 1.15831571885667
$ws.Range("C6").Value = 0.1647197264116755
$ws.Range("E6").Value = 0.09454989204672115
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 1.072737808054583
$ws.Range("H6").Value = 1.061332693993677
$ws.Range("L6").Value = 0.1855481797501639
$ws.Range("M6").Value = 0.2511297926447895
$ws.Range("N6").Value = 1.634954920013585

$ws.Range("B7").Value = 1.184066061757107
$ws.Range("C7").Value = 0.1692150174249321
$ws.Range("E7").Value = 0.09440839802856082
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 1.074546056503735
$ws.Range("H7").Value = 1.060243972777727
$ws.Range("L7").Value = 0.1862099183549049
$ws.Range("M7").Value = 0.2550087641154022
$ws.Range("N7").Value = 1.62839356068836

$ws.Range("B8").Value = 1.29899306187508
$ws.Range("C8").Value = 0.1888739359385454
$ws.Range("E8").Value = 0.09382082764458799
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 1.084397789859139
$ws.Range("H8").Value = 1.056793515364831
$ws.Range("L8").Value = 0.1893391156578872
$ws.Range("M8").Value = 0.2724301209722313
$ws.Range("N8").Value = 1.600871625736817

$ws.Range("B9").Value = 1.529025593851657
$ws.Range("C9").Value = 0.2267912409332666
$ws.Range("E9").Value = 0.09280035333885195
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 1.110505017616006
$ws.Range("H9").Value = 1.054911376582041
$ws.Range("L9").Value = 0.1962262588795269
$ws.Range("M9").Value = 0.307686974908961
$ws.Range("N9").Value = 1.552055236249474

$ws.Range("B10").Value = 1.700899153346484
$ws.Range("C10").Value = 0.2542709883035172
$ws.Range("E10").Value = 0.09213039975912096
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 1.133883841749707
$ws.Range("H10").Value = 1.056525420619892
$ws.Range("L10").Value = 0.2017444259903414
$ws.Range("M10").Value = 0.3342609514943362
$ws.Range("N10").Value = 1.519337280364566

$ws.Range("B11").Value = 1.779718992429821
$ws.Range("C11").Value = 0.2666920786650735
$ws.Range("E11").Value = 0.09184283379509472
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 1.145444165407753
$ws.Range("H11").Value = 1.057914760269796
$ws.Range("L11").Value = 0.2043544416632983
$ws.Range("M11").Value = 0.3464968894414611
$ws.Range("N11").Value = 1.505137685462277

$ws.Range("B12").Value = 1.809657319970086
$ws.Range("C12").Value = 0.2713842795383528
$ws.Range("E12").Value = 0.09173640493039237
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 1.149955798290222
$ws.Range("H12").Value = 1.058535399968264
$ws.Range("L12").Value = 0.2053571315477853
$ws.Range("M12").Value = 0.3511515273607557
$ws.Range("N12").Value = 1.499859176680873

$ws.Range("B13").Value = 1.803205522185237
$ws.Range("C13").Value = 0.2703742362145931
$ws.Range("E13").Value = 0.09175921672163967
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 1.148978163038322
$ws.Range("H13").Value = 1.058397524795453
$ws.Range("L13").Value = 0.205140547062328
$ws.Range("M13").Value = 0.3501481274105203
$ws.Range("N13").Value = 1.500991610576326

$ws.Range("B14").Value = 1.782180214450818
$ws.Range("C14").Value = 0.2670783369789547
$ws.Range("E14").Value = 0.09183402844407818
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 1.145812649358476
$ws.Range("H14").Value = 1.057963924342914
$ws.Range("L14").Value = 0.2044366463739067
$ws.Range("M14").Value = 0.3468794055057387
$ws.Range("N14").Value = 1.504701442628066

$ws.Range("B15").Value = 1.769313448176604
$ws.Range("C15").Value = 0.265058020810983
$ws.Range("E15").Value = 0.09188017373654933
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 1.143891158453073
$ws.Range("H15").Value = 1.057710651416045
$ws.Range("L15").Value = 0.2040073531370155
$ws.Range("M15").Value = 0.3448799729994576
$ws.Range("N15").Value = 1.506986665064534

$ws.Range("B16").Value = 1.695760829161372
$ws.Range("C16").Value = 0.2534576434396456
$ws.Range("E16").Value = 0.09214953835516493
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 1.133147045427563
$ws.Range("H16").Value = 1.056447833958856
$ws.Range("L16").Value = 0.2015758618406807
$ws.Range("M16").Value = 0.3334642658272102
$ws.Range("N16").Value = 1.520279038783134

$ws.Range("B17").Value = 1.650800867098724
$ws.Range("C17").Value = 0.2463208376727266
$ws.Range("E17").Value = 0.09231918488993651
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 1.126793514560177
$ws.Range("H17").Value = 1.055841153250981
$ws.Range("L17").Value = 0.2001097651498327
$ws.Range("M17").Value = 0.3264988063245227
$ws.Range("N17").Value = 1.528608811730023

$ws.Range("B18").Value = 1.62500074123534
$ws.Range("C18").Value = 0.2422084419477528
$ws.Range("E18").Value = 0.09241838038842021
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 1.123226156821431
$ws.Range("H18").Value = 1.055553848372057
$ws.Range("L18").Value = 0.199275897073008
$ws.Range("M18").Value = 0.3225063246412816
$ws.Range("N18").Value = 1.53346424179373

$ws.Range("B19").Value = 1.616275512804577
$ws.Range("C19").Value = 0.240814766917623
$ws.Range("E19").Value = 0.09245224461361579
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 1.122033228592898
$ws.Range("H19").Value = 1.05546714909886
$ws.Range("L19").Value = 0.1989951771783609
$ws.Range("M19").Value = 0.3211569213998331
$ws.Range("N19").Value = 1.535119255239112

$ws.Range("B20").Value = 1.655580758925112
$ws.Range("C20").Value = 0.2470813380643051
$ws.Range("E20").Value = 0.09230095818714301
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 1.127460845149614
$ws.Range("H20").Value = 1.055899353438235
$ws.Range("L20").Value = 0.200264861714686
$ws.Range("M20").Value = 0.3272388564717446
$ws.Range("N20").Value = 1.527715430126986

$ws.Range("B21").Value = 1.78835338944009
$ws.Range("C21").Value = 0.2680467313890347
$ws.Range("E21").Value = 0.09181198755414988
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 1.146738793435247
$ws.Range("H21").Value = 1.058088715272874
$ws.Range("L21").Value = 0.204643010165384
$ws.Range("M21").Value = 0.3478389344716604
$ws.Range("N21").Value = 1.503609098124248

$ws.Range("B22").Value = 1.875658127287579
$ws.Range("C22").Value = 0.2816824235933382
$ws.Range("E22").Value = 0.09150678819289015
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 1.16011948336353
$ws.Range("H22").Value = 1.060070717973957
$ws.Range("L22").Value = 0.2075879245118841
$ws.Range("M22").Value = 0.3614255649957556
$ws.Range("N22").Value = 1.488428860318612

$ws.Range("B23").Value = 1.829013486660926
$ws.Range("C23").Value = 0.2744108579365161
$ws.Range("E23").Value = 0.09166836619928054
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 1.152906140654807
$ws.Range("H23").Value = 1.058962349606276
$ws.Range("L23").Value = 0.2060085286441904
$ws.Range("M23").Value = 0.3541628524604548
$ws.Range("N23").Value = 1.496478186113838

$ws.Range("B24").Value = 1.653419620445277
$ws.Range("C24").Value = 0.246737544970415
$ws.Range("E24").Value = 0.0923091932977258
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 1.127158879139444
$ws.Range("H24").Value = 1.055872849639627
$ws.Range("L24").Value = 0.200194714491019
$ws.Range("M24").Value = 0.3269042422527093
$ws.Range("N24").Value = 1.528119120720211

$ws.Range("B25").Value = 1.466293937935518
$ws.Range("C25").Value = 0.2166006050215401
$ws.Range("E25").Value = 0.09306236894773923
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 1.102709467567777
$ws.Range("H25").Value = 1.054895674789762
$ws.Range("L25").Value = 0.1942826793285946
$ws.Range("M25").Value = 0.298031568893137
$ws.Range("N25").Value = 1.564709051518709
